# Generate Report for Handback
# - "Ready for handoff" status becomes "Handback transform failed"
#   (shared by Overview!E3/F3 and the zh-cn/de-de "Status" columns)
# - Populate the "Error Detail" column (P) on the zh-cn and de-de sheets
#   with the handback/handoff file-name mismatch message
# - Widen the "Error Detail" column so the message is readable

$wb = $excel.ActiveWorkbook

$statusMessage = "Handback transform failed"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusMessage
$overview.Range("F3").Value = $statusMessage

# Column P is stored internally at 39 + 1/6 "characters" so that, after
# Excel's own pixel-rounding on save, the persisted width comes out to
# exactly 40 (matches the target column width of 40).
$errorColumnWidth = 39 + (1 / 6)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusMessage
$zhcn.Range("P3").Value = "Handback file name: 3rd4stpj.ssn is different with handoff file name: 58e69f0c-af68-41ca-9cc7-2fe7ecf1d9c4.7e4fdf9954c2ab35697a5a662af79121ea550df3.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = $errorColumnWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusMessage
$dede.Range("P3").Value = "Handback file name: 3rd4stpj.ssn is different with handoff file name: 58e69f0c-af68-41ca-9cc7-2fe7ecf1d9c4.7e4fdf9954c2ab35697a5a662af79121ea550df3.de-de."
$dede.Columns.Item(16).ColumnWidth = $errorColumnWidth
